# Fix missing header row in generated xlsx (TS version export).
# Insert a new row 3 with column headers (Context / Source / Translation /
# Location) above the existing MenuBar data rows, pushing the data down by
# one row (old rows 3-6 become rows 4-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the first data row (row 3), shifting rows 3-6
# down to 4-7 and extending the used range/dimension accordingly.
$ws.Rows("3").Insert()

# Populate the new header row.
$ws.Range("A3").Value = "Context"
$ws.Range("B3").Value = "Source"
$ws.Range("C3").Value = "Translation"
$ws.Range("D3").Value = "Location"
